# Apply weekly price/volume refresh to cryptos sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.553.93'
$ws.Range("E2").Value = '  +1.05%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.392.87'
$ws.Range("E3").Value = '  +0.15%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.60'
$ws.Range("E5").Value = '  +0.65%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.09'
$ws.Range("E6").Value = '  -0.77%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("E8").Value = '  -0.32%  '
$ws.Range("E9").Value = '  +2.52%  '
$ws.Range("E10").Value = '  -0.76%  '
$ws.Range("E11").Value = '  -1.88%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.974.48'
$ws.Range("E12").Value = '  +0.24%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.125'
$ws.Range("E13").Value = '  +0.01%  '
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.42'
$ws.Range("E14").Value = '  +0.82%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.398.56'
$ws.Range("E15").Value = '  +0.33%  '
$ws.Range("E16").Value = '  -0.18%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.589.37'
$ws.Range("E17").Value = '  +1.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.15'
$ws.Range("E18").Value = '  -0.29%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.69'
$ws.Range("E19").Value = '  -1.00%  '
$ws.Range("E20").Value = '  +0.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '391.05'
$ws.Range("E21").Value = '  +1.76%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '75.35'
$ws.Range("E22").Value = '  +1.22%  '
$ws.Range("E23").Value = '  -0.62%  '
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("E25").Value = '  -4.19%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.192'
$ws.Range("E26").Value = '  +7.03%  '
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.26'
$ws.Range("E28").Value = '  -1.81%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.02'
$ws.Range("E29").Value = '  +0.61%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.14'
$ws.Range("E30").Value = '  +0.06%  '
$ws.Range("E31").Value = '  -0.05%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.38'
$ws.Range("E32").Value = '  -3.43%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.34'
$ws.Range("E33").Value = '  -0.72%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.93'
$ws.Range("E34").Value = '  -0.79%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '168.06'
$ws.Range("E35").Value = '  +0.17%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.06'
$ws.Range("E36").Value = '  +1.48%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.429.60'
$ws.Range("E37").Value = '  +0.40%  '
$ws.Range("E38").Value = '  -1.13%  '
$ws.Range("E39").Value = '  -0.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '25.98'
$ws.Range("E40").Value = '  -5.77%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.780'
$ws.Range("E41").Value = '  +0.09%  '
$ws.Range("E42").Value = '  -0.09%  '
$ws.Range("E43").Value = '  -0.83%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.15'
$ws.Range("E44").Value = '  +1.16%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.476.12'
$ws.Range("E45").Value = '  +0.10%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '23.02'
$ws.Range("E46").Value = '  +0.07%  '
$ws.Range("E47").Value = '  -2.28%  '
$ws.Range("E48").Value = '  +0.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0263'
$ws.Range("E49").Value = '  -1.58%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.06'
$ws.Range("E50").Value = '  -1.47%  '
$ws.Range("E51").Value = '  -1.70%  '
